{"js": "// Minor wording fixes to the \"Change Requests\" table:\n//  1. \"Right Ascension and Declination\" -> \"right ascension and declination\"\n//  2. \"user guide on the Admin user\"    -> \"user's guide on the admin user\"\n\nconst body = context.document.body;\n\n// --- Change 1: de-capitalize \"Right Ascension and Declination\" ---------\nconst r1 = body.search(\"Right Ascension and Declination\", { matchCase: true });\nr1.load(\"text\");\nawait context.sync();\n\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\"right ascension and declination\", Word.InsertLocation.replace);\n}\n\n// --- Change 2: \"user guide on the Admin user\" -> \"user's guide on the admin user\" ---\nconst r2 = body.search(\"user guide on the Admin user\", { matchCase: true });\nr2.load(\"text\");\nawait context.sync();\n\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\"user\\u2019s guide on the admin user\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Minor wording fixes to the \"Change Requests\" table:\n#  1. \"Right Ascension and Declination\" -> \"right ascension and declination\"\n#  2. \"user guide on the Admin user\"    -> \"user's guide on the admin user\"\n\n$d = $word.ActiveDocument\n\n# --- Change 1: de-capitalize \"Right Ascension and Declination\" ---------\n$find1 = $d.Content.Find\n$found1 = $find1.Execute(\n    \"Right Ascension and Declination\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"right ascension and declination\",\n    2\n)\n\n# --- Change 2: \"user guide on the Admin user\" -> \"user's guide on the admin user\" ---\n$find2 = $d.Content.Find\n$apostropheS = [char]0x2019 + \"s\"\n$found2 = $find2.Execute(\n    \"user guide on the Admin user\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"user\" + $apostropheS + \" guide on the admin user\",\n    2\n)\n"}
